# "Updated example to have a scale column"
# Adds a new "Scale" column (D) to the position-file example sheet:
#   - D9  gets the header label "Scale" (same text already used in A3)
#   - D10:D64 each get a scale factor of 1, one per montage position row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header
$ws.Range("D9").Value = "Scale"

# Scale factor of 1 for every data row below the header
$ws.Range("D10:D64").Value = 1

# Leave the selection where the edit was made, like the authored workbook
$ws.Range("D11:D64").Select() | Out-Null
